$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$headerRange = $ws.Range("A1:U1")
$scratch = $ws.Range("A100:U100")
$headerRange.Copy()
$scratch.PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)
$headerRange.ClearFormats()
$range = $ws.Range("A1:U55")
$lo = $ws.ListObjects.Add([Microsoft.Office.Interop.Excel.XlListObjectSourceType]::xlSrcRange, $range, $null, [Microsoft.Office.Interop.Excel.XlYesNoGuess]::xlYes)
Write-Host ("Created table: " + $lo.Name)
$scratch.Copy()
$headerRange.PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)
$scratch.ClearFormats()
$scratch.ClearContents()
